# Punto de Control 2
# Fill in inventory rows 3-6 with product/quantity data and normalize
# the cell formatting of the data range so it matches the header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "tecnologia"
$ws.Range("B3").Value = "3"
$ws.Range("A4").Value = "carro"
$ws.Range("B4").Value = "2"
$ws.Range("A5").Value = "guadañadora"
$ws.Range("B5").Value = "3"
$ws.Range("A6").Value = "taladro"
$ws.Range("B6").Value = "4"

# Rows 3-12 previously carried a distinct (applyFont=false) text style;
# re-apply the text number format across the whole table so the cells
# consolidate onto the same style as the header/first data rows.
$ws.Range("A1:B12").NumberFormat = "@"

$ws.Range("B6").Select()
